# Swap the order of the two "Shopping Behavior" slides (positions 2 and 3):
# the "An exploratory study of grocery shopping stressors" slide now comes
# before the "Shopping in physical stores" slide.
$p = $ppt.ActivePresentation

$stressorsSlide = $p.Slides.Item(3)
$stressorsSlide.MoveTo(2)
